$d = $word.ActiveDocument

# The bibliography entry for BATALHA is followed by: a blank paragraph,
# a "Ver no Jupiter..." paragraph, and a "(c) 2020 ..." paragraph. Those
# trailing three paragraphs (the footer scraped from the course-catalog
# website) need to be removed, leaving the BATALHA paragraph followed
# directly by the single blank paragraph that used to sit right before
# the final page-break paragraph.

$range = $d.Range(0, 0)
$range.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx")
$jupiterPara = $range.Paragraphs(1)

# Paragraph right before "Ver no Jupiter..." (the blank separator) and the
# paragraph right after it (the copyright/footer line) bound the block that
# must disappear.
$blankBefore = $jupiterPara.Previous(1)
$copyrightPara = $jupiterPara.Next(1)

$deleteRange = $d.Range($blankBefore.Range.Start, $copyrightPara.Range.End)
$deleteRange.Delete()
